# Auto-generated edit script for cryptos.xlsx update
# Commit: "Updated cryptos list on Wed Nov 20 23:44:18 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain decimal number need to be
# force-formatted as Text first, otherwise Excel COM auto-converts the assigned
# string into a numeric Value and trailing/insignificant zeros would be lost
# (e.g. "8.90" -> 8.9, "0.000270" -> 0.00027), which would not match the workbook
# convention of storing these Price values as literal text.
$textCells = @('D5', 'D6', 'D15', 'D16', 'D20', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D44', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '94.124.43'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '3.079.38'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '236.56'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '608.36'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +9.11%  '
$ws.Range('D11').Value = '3.073.91'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('D13').Value = '93.786.32'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('E14').Value = '  -3.08%  '
$ws.Range('D15').Value = '33.82'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '5.31'
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('D17').Value = '3.644.91'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').Value = '3.076.24'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('E19').Value = '  -4.25%  '
$ws.Range('D20').Value = '14.31'
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('D22').Value = '441.67'
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('D23').Value = '8.86'
$ws.Range('E23').Value = '  -5.64%  '
$ws.Range('D24').Value = '0.0000190'
$ws.Range('E24').Value = '  -3.89%  '
$ws.Range('D25').Value = '8.38'
$ws.Range('E25').Value = '  +6.24%  '
$ws.Range('E26').Value = '  -4.56%  '
$ws.Range('D27').Value = '84.49'
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('D28').Value = '11.91'
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('D29').Value = '3.230.97'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +7.03%  '
$ws.Range('E32').Value = '  +5.14%  '
$ws.Range('E33').Value = '  -7.54%  '
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('D35').Value = '8.90'
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('D36').Value = '7.47'
$ws.Range('E36').Value = '  -5.24%  '
$ws.Range('E37').Value = '  -4.03%  '
$ws.Range('D38').Value = '25.51'
$ws.Range('E38').Value = '  -2.63%  '
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').Value = '1.88'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '485.04'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').Value = '3.84'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').Value = '1.25'
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('E46').Value = '  -6.78%  '
$ws.Range('D47').Value = '161.10'
$ws.Range('E47').Value = '  -0.57%  '
$ws.Range('D48').Value = '0.676'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').Value = '1.83'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '0.000270'
$ws.Range('E50').Value = '  +8.88%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '43.55'
$ws.Range('E51').Value = '  -0.94%  '
